# Applies: add Sheet3 (Sales Report) and Sheet4 (Target Goals) with data,
# formulas and formatting; update Sheet1/Sheet2 selections; make Sheet3 active.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add two new worksheets at the end: Sheet3, Sheet4
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet2)

# ---------------------------------------------------------------------
# Sheet3: Sales Report data
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = "Salesperson"
$ws3.Range("B1").Value = "Region"
$ws3.Range("C1").Value = "Product"
$ws3.Range("D1").Value = "Sales Amount"
$ws3.Range("E1").Value = "Target"
$ws3.Range("F1").Value = "TARGET GOALS"

$salesData = @(
    @("Neha",   "North", "Laptop", 45000, 40000),
    @("Rajesh", "South", "Mobile", 38000, 35000),
    @("Priya",  "East",  "Tablet", 25000, 30000),
    @("Suresh", "West",  "Laptop", 60000, 50000),
    @("Anita",  "North", "Mobile", 42000, 40000),
    @("Ramesh", "East",  "Laptop", 55000, 45000),
    @("Kavita", "South", "Tablet", 28000, 30000)
)

$r = 2
foreach ($row in $salesData) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# F2:F8 - Nested/Simple IF Achieved vs Not Achieved
$ws3.Range("F2:F8").Formula = '=IF(D2>=E2,"ACHIEVED","NOT ACHIEVED")'

# Question labels + answers block
$ws3.Range("D12").Value = "2.COUNT IF"
$ws3.Range("E12").Formula = '=COUNTIF(B2:B8,"NORTH")'

$ws3.Range("D13").Value = "3.COUNTIFS"
$ws3.Range("E13").Formula = '=COUNTIFS(B2:B8,"EAST",C2:C8,"LAPTOP")'

$ws3.Range("D14").Value = "4.SUMIF"
$ws3.Range("E14").Formula = '=SUMIF(C2:C8,"MOBILE",D2:D8)'

$ws3.Range("D15").Value = "5.SUMIFS"
$ws3.Range("E15").Formula = '=SUMIFS(D2:D8,B2:B8,"SOUTH",C2:C8,"TABLET")'

$ws3.Range("D17").Value = "6.VLOOKUP"
$ws3.Range("E17").Formula = '=VLOOKUP("SURESH",A2:E8,5,FALSE)'

$ws3.Range("D19").Value = "7.MEAN"
$ws3.Range("E19").Formula = "=AVERAGE(D2:D8)"

$ws3.Range("D20").Value = "MEDIAN"
$ws3.Range("E20").Formula = "=MEDIAN(D2:D8)"

$ws3.Range("D21").Value = "MODE"
$ws3.Range("E21").Formula = "=_xlfn.MODE.SNGL(D2:D8)"

# Column widths on Sheet3
$ws3.Columns.Item(1).ColumnWidth = 11.5546875
$ws3.Columns.Item(4).ColumnWidth = 12.44140625
$ws3.Columns.Item(6).ColumnWidth = 15.33203125

# Formatting
$ws3.Range("F1:F8").Font.Color = 255       # red font
$ws3.Range("F1:F8").Interior.Color = 65535 # yellow fill

$ws3.Range("D12:E15").Interior.Color = 5296274   # light green FF92D050

$ws3.Range("D17:E17").Interior.Color = 255        # red fill FFFF0000

$ws3.Range("D19:E21").Interior.Color = 5287936    # green FF00B050

$ws3.Range("D19:E21").Select()

# ---------------------------------------------------------------------
# Sheet4: Sales Report Questions
# ---------------------------------------------------------------------
$questions = @(
    "Sales Report Questions",
    "1. Nested IF: Display 'Achieved' if Sales >= Target, else 'Not Achieved'.",
    "2. COUNTIF: Count salespersons in 'North' region.",
    "3. COUNTIFS: Count salespersons who sold Laptop in East region.",
    "4. SUMIF: Total sales of Mobile category.",
    "5. SUMIFS: Total sales in South region for Tablet.",
    "6. VLOOKUP: Retrieve Target for 'Suresh'.",
    "7. Mean, Median, Mode of Sales Amount."
)

$r = 1
foreach ($q in $questions) {
    $ws4.Cells.Item($r, 1).Value = $q
    $r = $r + 1
}

$ws4.Range("A1").Font.Bold = $true
$ws4.Range("A1").Borders.LineStyle = 1
$ws4.Range("A1").HorizontalAlignment = -4108
$ws4.Range("A1").VerticalAlignment = -4160

$ws4.Columns.Item(1).ColumnWidth = 11.109375
$ws4.Columns.Item(2).ColumnWidth = 11.33203125
$ws4.Columns.Item(3).ColumnWidth = 11.33203125
$ws4.Columns.Item(4).ColumnWidth = 10.77734375
$ws4.Columns.Item(5).ColumnWidth = 10.5546875
$ws4.Columns.Item(6).ColumnWidth = 10.5546875

$ws4.Range("D11").Select()

# ---------------------------------------------------------------------
# Other workbook-level changes
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("D26").Select()

# Make Sheet3 the active/selected tab (last action wins for tabSelected)
$ws3.Activate()
